$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 91

$ws.Cells.Item($row, 1).Value = "2025/12/06 14:00"

for ($col = 2; $col -le 7; $col++) {
    $ws.Cells.Item($row, $col).Value = "-"
}
